# Update the "想去人数" (want-to-go count) figures for two 2024-10 events.
# These counts are duplicated on both the "展览" sheet and the combined
# "全部类型" sheet, so both need the same update:
#   row 4 (青田·未闻展名国漫嘉年华): F4  79 -> 80
#   row 5 (缙云·星辰动漫游戏展嘉年华): F5  19 -> 21

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 80
    $ws.Range("F5").Value = 21
}
